# Update to Blank correction & Mean correction
#
# 1. On the "Config" sheet, the IntraBatchMode row ("Median","Linear","Spline")
#    is renamed to ("Mean","Linear","Spline") and its explanation text is
#    updated to talk about "Mean" instead of "Median".
# 2. The active sheet/tab switches from "Clean & Explore" to "Config", with the
#    Config sheet's selection landing on C8.

$wb = $excel.ActiveWorkbook

$wsConfig = $wb.Worksheets.Item("Config")

# --- content changes -------------------------------------------------
$wsConfig.Range("B8").Value = """Mean"",""Linear"",""Spline"""
$wsConfig.Range("C8").Value = "Three correction modes. ""Spline"" is the default QCRSC algorithm that requires optimisation of the smoothing parameter. ""Linear"" is a simple (interpolated) linear regression based on the QC values & as such requires no smoothing optimisation. ""Mean"" equalises the QC mean across batches & ignores within batch systematic change."

# --- view / selection changes -----------------------------------------
# Make the Config sheet the active tab and leave the selection on C8, matching
# the workbook's new activeTab + sheetView selection state.
$wsConfig.Activate()
$wsConfig.Range("C8").Select()
